$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $r = $ws.Cells.Item($row, $col)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

function Set-Cell($row, $col, $value) {
    $r = $ws.Cells.Item($row, $col)
    $r.Value = $value
}

Set-TextCell 2 4 '25.995.36'
Set-TextCell 2 5 '  -0.28%  '

Set-TextCell 3 4 '1.743.82'
Set-TextCell 3 5 '  -0.22%  '

Set-TextCell 4 5 '  -0.09%  '

Set-TextCell 5 4 '250.76'
Set-TextCell 5 5 '  +7.46%  '

Set-TextCell 6 5 '  -0.15%  '

Set-TextCell 7 4 '0.5154'
Set-TextCell 7 5 '  -2.00%  '

Set-Cell 8 2 'OKB'
Set-Cell 8 3 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 8 4 '41.16'
Set-TextCell 8 5 '  +2.57%  '

Set-Cell 9 2 'Cardano'
Set-Cell 9 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 9 4 '0.2762'
Set-TextCell 9 5 '  -0.12%  '

Set-Cell 10 2 'Dogecoin'
Set-Cell 10 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 10 4 '0.06197'
Set-TextCell 10 5 '  +0.17%  '

Set-Cell 11 2 'WrappedEther'
Set-Cell 11 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 11 4 '1.743.00'
Set-TextCell 11 5 '  -0.26%  '

Set-Cell 12 2 'TRON'
Set-Cell 12 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 12 4 '0.07239'
Set-TextCell 12 5 '  +0.38%  '

Set-Cell 13 2 'Solana'
Set-Cell 13 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 13 4 '15.21'
Set-TextCell 13 5 '  -0.29%  '

Set-Cell 14 2 'Polygon'
Set-Cell 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 14 4 '0.6524'
Set-TextCell 14 5 '  +2.10%  '

Set-Cell 15 2 'Polkadot'
Set-Cell 15 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 15 4 '4.635'
Set-TextCell 15 5 '  +1.09%  '

Set-Cell 16 2 'Litecoin'
Set-Cell 16 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 16 4 '77.92'
Set-TextCell 16 5 '  -0.39%  '

Set-Cell 17 2 'Dai'
Set-Cell 17 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 17 4 '0.9997'
Set-TextCell 17 5 '  -0.16%  '

Set-Cell 18 2 'BinanceUSD'
Set-Cell 18 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 18 4 '0.9996'
Set-TextCell 18 5 '  -0.08%  '

Set-Cell 19 2 'WrappedBTC'
Set-Cell 19 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 19 4 '26.018.95'
Set-TextCell 19 5 '  +0.10%  '

Set-Cell 20 2 'Avalanche'
Set-Cell 20 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 20 4 '11.87'
Set-TextCell 20 5 '  +2.73%  '

Set-Cell 21 2 'ShibaInu'
Set-Cell 21 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 21 4 '0.000006819'
Set-TextCell 21 5 '  +1.84%  '

Set-Cell 22 2 'WrappedliquidstakedEther2.0'
Set-Cell 22 3 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextCell 22 4 '1.966.23'
Set-TextCell 22 5 '  -0.50%  '

Set-Cell 23 2 'Uniswap'
Set-Cell 23 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 23 4 '4.305'
Set-TextCell 23 5 '  -0.42%  '

Set-Cell 24 2 'Cosmos'
Set-Cell 24 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 24 4 '8.682'
Set-TextCell 24 5 '  -1.39%  '

Set-Cell 25 2 'Chainlink'
Set-Cell 25 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 25 4 '5.393'
Set-TextCell 25 5 '  +3.93%  '

Set-Cell 26 2 'Monero'
Set-Cell 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 26 4 '136.09'
Set-TextCell 26 5 '  -2.48%  '

Set-Cell 27 2 'Toncoin'
Set-Cell 27 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 27 4 '1.510'
Set-TextCell 27 5 '  -0.89%  '

Set-Cell 28 2 'EthereumClassic'
Set-Cell 28 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 28 4 '15.29'
Set-TextCell 28 5 '  +0.18%  '

Set-Cell 29 2 'LidoDAOToken'
Set-Cell 29 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 29 4 '1.792'
Set-TextCell 29 5 '  -0.77%  '

Set-Cell 30 2 'BitcoinCash'
Set-Cell 30 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 30 4 '105.98'
Set-TextCell 30 5 '  +1.85%  '

Set-Cell 31 2 'InternetComputer(DFINITY)'
Set-Cell 31 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 31 4 '3.963'
Set-TextCell 31 5 '  +5.15%  '

Set-Cell 32 2 'Stellar'
Set-Cell 32 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 32 4 '0.08252'
Set-TextCell 32 5 '  -0.37%  '

Set-Cell 33 2 'Filecoin'
Set-Cell 33 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 33 4 '3.675'
Set-TextCell 33 5 '  +0.59%  '

Set-Cell 34 2 'Hedera'
Set-Cell 34 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 34 4 '0.04687'
Set-TextCell 34 5 '  +3.71%  '

Set-Cell 35 2 'HuobiToken'
Set-Cell 35 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 35 4 '2.658'
Set-TextCell 35 5 '  +0.71%  '

Set-Cell 36 2 'ARBITRUM'
Set-Cell 36 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 36 4 '1.002'
Set-TextCell 36 5 '  +0.56%  '

Set-Cell 37 2 'ImmutableX'
Set-Cell 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 37 4 '0.6266'
Set-TextCell 37 5 '  -0.20%  '

Set-Cell 38 2 'MXToken'
Set-Cell 38 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 38 4 '2.730'
Set-TextCell 38 5 '  +1.08%  '

Set-Cell 39 2 'VeChain'
Set-Cell 39 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 39 4 '0.01615'
Set-TextCell 39 5 '  +1.65%  '

Set-Cell 40 2 'RenderToken'
Set-Cell 40 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 40 4 '1.930'
Set-TextCell 40 5 '  +0.54%  '

Set-Cell 41 2 'PaxDollar'
Set-Cell 41 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 41 4 '0.9994'
Set-TextCell 41 5 '  -0.20%  '

Set-Cell 42 2 'Quant'
Set-Cell 42 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 42 4 '100.76'
Set-TextCell 42 5 '  +3.25%  '

Set-TextCell 43 4 '0.7627'
Set-TextCell 43 5 '  +4.15%  '

Set-Cell 44 2 'TheSandbox'
Set-Cell 44 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 44 4 '0.3889'
Set-TextCell 44 5 '  +0.07%  '

Set-Cell 45 2 'FraxShare'
Set-Cell 45 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 45 4 '5.025'
Set-TextCell 45 5 '  -0.14%  '

Set-Cell 46 2 'Aptos'
Set-Cell 46 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 46 4 '6.373'
Set-TextCell 46 5 '  +1.35%  '

Set-Cell 47 2 'Algorand'
Set-Cell 47 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 47 4 '0.1134'
Set-TextCell 47 5 '  -0.36%  '

Set-Cell 48 2 'Aave'
Set-Cell 48 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 48 4 '55.65'
Set-TextCell 48 5 '  +3.31%  '

Set-Cell 49 2 'Cronos'
Set-Cell 49 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 49 4 '0.05231'
Set-TextCell 49 5 '  -2.05%  '

Set-Cell 50 2 'Elrond'
Set-Cell 50 3 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 50 4 '30.83'
Set-TextCell 50 5 '  +1.11%  '

Set-Cell 51 2 'Decentraland'
Set-Cell 51 3 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 51 4 '0.3449'
Set-TextCell 51 5 '  +0.05%  '
